# Update LR-pair TPM-derived metrics (Icam4-Itga2b) with newly computed
# TPM-based values across rows 2-17 (commit: "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.389221
$ws.Range("H2").Value = 4.167663
$ws.Range("I2").Value = 0.2910270461264192
$ws.Range("J2").Value = 0.2910270461264192
$ws.Range("M2").Value = 1.646992
$ws.Range("N2").Value = 4.940976
$ws.Range("O2").Value = 0.2071783517404009
$ws.Range("P2").Value = 0.2071783517404009
$ws.Range("Q2").Value = 2.288035873232
$ws.Range("R2").Value = 20.592322859088
$ws.Range("S2").Value = 0.06029450372834916
$ws.Range("T2").Value = 0.06029450372834916

# Row 3
$ws.Range("G3").Value = 1.389221
$ws.Range("H3").Value = 4.167663
$ws.Range("I3").Value = 0.2910270461264192
$ws.Range("J3").Value = 0.2910270461264192
$ws.Range("O3").Value = 0.4685125322965616
$ws.Range("P3").Value = 0.4685125322965616
$ws.Range("Q3").Value = 5.174157782163
$ws.Range("R3").Value = 46.567420039467
$ws.Range("S3").Value = 0.1363498183474769
$ws.Range("T3").Value = 0.1363498183474769

# Row 4
$ws.Range("G4").Value = 1.389221
$ws.Range("H4").Value = 4.167663
$ws.Range("I4").Value = 0.2910270461264192
$ws.Range("J4").Value = 0.2910270461264192
$ws.Range("M4").Value = 2.284352333333333
$ws.Range("N4").Value = 6.853057
$ws.Range("O4").Value = 0.2873531572796583
$ws.Range("P4").Value = 0.2873531572796583
$ws.Range("Q4").Value = 3.173470232865666
$ws.Range("R4").Value = 28.561232095791
$ws.Range("S4").Value = 0.08362754055819929
$ws.Range("T4").Value = 0.08362754055819929

# Row 5
$ws.Range("G5").Value = 1.389221
$ws.Range("H5").Value = 4.167663
$ws.Range("I5").Value = 0.2910270461264192
$ws.Range("J5").Value = 0.2910270461264192
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2937863333333333
$ws.Range("N5").Value = 0.881359
$ws.Range("O5").Value = 0.03695595868337916
$ws.Range("P5").Value = 0.03695595868337916
$ws.Range("Q5").Value = 0.4081341437796667
$ws.Range("R5").Value = 3.673207294017
$ws.Range("S5").Value = 0.01075518349239383
$ws.Range("T5").Value = 0.01075518349239383

# Row 6
$ws.Range("I6").Value = 0.461328155686921
$ws.Range("J6").Value = 0.4613281556869209
$ws.Range("M6").Value = 1.646992
$ws.Range("N6").Value = 4.940976
$ws.Range("O6").Value = 0.2071783517404009
$ws.Range("P6").Value = 0.2071783517404009
$ws.Range("Q6").Value = 3.626932216757334
$ws.Range("R6").Value = 32.642389950816
$ws.Range("S6").Value = 0.09557720690665536
$ws.Range("T6").Value = 0.09557720690665535

# Row 7
$ws.Range("I7").Value = 0.461328155686921
$ws.Range("J7").Value = 0.4613281556869209
$ws.Range("O7").Value = 0.4685125322965616
$ws.Range("P7").Value = 0.4685125322965616
$ws.Range("S7").Value = 0.2161380224405818
$ws.Range("T7").Value = 0.2161380224405817

# Row 8
$ws.Range("I8").Value = 0.461328155686921
$ws.Range("J8").Value = 0.4613281556869209
$ws.Range("M8").Value = 2.284352333333333
$ws.Range("N8").Value = 6.853057
$ws.Range("O8").Value = 0.2873531572796583
$ws.Range("P8").Value = 0.2873531572796583
$ws.Range("Q8").Value = 5.030498674062445
$ws.Range("R8").Value = 45.27448806656201
$ws.Range("S8").Value = 0.1325641020786385
$ws.Range("T8").Value = 0.1325641020786385

# Row 9
$ws.Range("I9").Value = 0.461328155686921
$ws.Range("J9").Value = 0.4613281556869209
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2937863333333333
$ws.Range("N9").Value = 0.881359
$ws.Range("O9").Value = 0.03695595868337916
$ws.Range("P9").Value = 0.03695595868337916
$ws.Range("Q9").Value = 0.6469631408104445
$ws.Range("R9").Value = 5.822668267294001
$ws.Range("S9").Value = 0.01704882426104536
$ws.Range("T9").Value = 0.01704882426104536

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1506176666666667
$ws.Range("H10").Value = 0.451853
$ws.Range("I10").Value = 0.03155280162368235
$ws.Range("J10").Value = 0.03155280162368235
$ws.Range("M10").Value = 1.646992
$ws.Range("N10").Value = 4.940976
$ws.Range("O10").Value = 0.2071783517404009
$ws.Range("P10").Value = 0.2071783517404009
$ws.Range("Q10").Value = 0.2480660920586667
$ws.Range("R10").Value = 2.232594828528
$ws.Range("S10").Value = 0.006537057433186357
$ws.Range("T10").Value = 0.006537057433186357

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1506176666666667
$ws.Range("H11").Value = 0.451853
$ws.Range("I11").Value = 0.03155280162368235
$ws.Range("J11").Value = 0.03155280162368235
$ws.Range("O11").Value = 0.4685125322965616
$ws.Range("P11").Value = 0.4685125322965616
$ws.Range("Q11").Value = 0.560975951353
$ws.Range("R11").Value = 5.048783562177
$ws.Range("S11").Value = 0.01478288298976248
$ws.Range("T11").Value = 0.01478288298976248

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.1506176666666667
$ws.Range("H12").Value = 0.451853
$ws.Range("I12").Value = 0.03155280162368235
$ws.Range("J12").Value = 0.03155280162368235
$ws.Range("M12").Value = 2.284352333333333
$ws.Range("N12").Value = 6.853057
$ws.Range("O12").Value = 0.2873531572796583
$ws.Range("P12").Value = 0.2873531572796583
$ws.Range("Q12").Value = 0.3440638182912222
$ws.Range("R12").Value = 3.096574364621
$ws.Range("S12").Value = 0.009066797167583853
$ws.Range("T12").Value = 0.009066797167583853

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.1506176666666667
$ws.Range("H13").Value = 0.451853
$ws.Range("I13").Value = 0.03155280162368235
$ws.Range("J13").Value = 0.03155280162368235
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.2937863333333333
$ws.Range("N13").Value = 0.881359
$ws.Range("O13").Value = 0.03695595868337916
$ws.Range("P13").Value = 0.03695595868337916
$ws.Range("Q13").Value = 0.04424941202522222
$ws.Range("R13").Value = 0.398244708227
$ws.Range("S13").Value = 0.001166064033149664
$ws.Range("T13").Value = 0.001166064033149664

# Row 14
$ws.Range("G14").Value = 1.031517666666667
$ws.Range("H14").Value = 3.094553
$ws.Range("I14").Value = 0.2160919965629775
$ws.Range("J14").Value = 0.2160919965629775
$ws.Range("M14").Value = 1.646992
$ws.Range("N14").Value = 4.940976
$ws.Range("O14").Value = 0.2071783517404009
$ws.Range("P14").Value = 0.2071783517404009
$ws.Range("Q14").Value = 1.698901344858667
$ws.Range("R14").Value = 15.290112103728
$ws.Range("S14").Value = 0.04476958367221007
$ws.Range("T14").Value = 0.04476958367221007

# Row 15
$ws.Range("G15").Value = 1.031517666666667
$ws.Range("H15").Value = 3.094553
$ws.Range("I15").Value = 0.2160919965629775
$ws.Range("J15").Value = 0.2160919965629775
$ws.Range("O15").Value = 0.4685125322965616
$ws.Range("P15").Value = 0.4685125322965616
$ws.Range("Q15").Value = 3.841890644053
$ws.Range("R15").Value = 34.577015796477
$ws.Range("S15").Value = 0.1012418085187405
$ws.Range("T15").Value = 0.1012418085187405

# Row 16
$ws.Range("G16").Value = 1.031517666666667
$ws.Range("H16").Value = 3.094553
$ws.Range("I16").Value = 0.2160919965629775
$ws.Range("J16").Value = 0.2160919965629775
$ws.Range("M16").Value = 2.284352333333333
$ws.Range("N16").Value = 6.853057
$ws.Range("O16").Value = 0.2873531572796583
$ws.Range("P16").Value = 0.2873531572796583
$ws.Range("Q16").Value = 2.356349788724555
$ws.Range("R16").Value = 21.207148098521
$ws.Range("S16").Value = 0.06209471747523666
$ws.Range("T16").Value = 0.06209471747523666

# Row 17
$ws.Range("G17").Value = 1.031517666666667
$ws.Range("H17").Value = 3.094553
$ws.Range("I17").Value = 0.2160919965629775
$ws.Range("J17").Value = 0.2160919965629775
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.2937863333333333
$ws.Range("N17").Value = 0.881359
$ws.Range("O17").Value = 0.03695595868337916
$ws.Range("P17").Value = 0.03695595868337916
$ws.Range("Q17").Value = 0.3030457930585555
$ws.Range("R17").Value = 2.727412137527
$ws.Range("S17").Value = 0.007985886896790308
$ws.Range("T17").Value = 0.00798588689679031
